$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the country names for rows 156 and 157 (Uganda moves above Mozambique
# in the ranking, taking Mozambique's old position and vice versa).
$ws.Range("A156").Value = "Uganda"
$ws.Range("A157").Value = "Mozambique"

# Update Estados Unidos (row 4) stats
$ws.Range("B4").Value = 985060
$ws.Range("C4").Value = 24409
$ws.Range("E4").Value = 810968
$ws.Range("G4").Value = 1101
$ws.Range("H4").Value = 55357

# Update Alemania (row 8) stats
$ws.Range("B8").Value = 157495
$ws.Range("C8").Value = 982
$ws.Range("E8").Value = 39551
$ws.Range("G8").Value = 67
$ws.Range("H8").Value = 5944

# Update Tunez (row 90) stats
$ws.Range("B90").Value = 949
$ws.Range("C90").Value = 10
$ws.Range("D90").Value = 216
$ws.Range("E90").Value = 695
$ws.Range("F90").Value = 20

# Update row 156 stats (now Uganda)
$ws.Range("B156").Value = 79
$ws.Range("C156").Value = 4
$ws.Range("D156").Value = 46
$ws.Range("E156").Value = 33

# Update row 157 stats (now Mozambique)
$ws.Range("B157").Value = 76
$ws.Range("C157").Value = 6
$ws.Range("D157").Value = 12
$ws.Range("E157").Value = 64
